$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1: clear existing format, then give it a top+bottom thin border (matches borderId=4)
$ws1.Range("C1").ClearFormats()
$c1b = $ws1.Range("C1").Borders
$c1b.Item(8).ColorIndex = -4105
$c1b.Item(8).LineStyle = 1
$c1b.Item(9).ColorIndex = -4105
$c1b.Item(9).LineStyle = 1

# D1: clear existing format, then give it a top+right+bottom thin border (matches borderId=5)
$ws1.Range("D1").ClearFormats()
$d1b = $ws1.Range("D1").Borders
$d1b.Item(8).ColorIndex = -4105
$d1b.Item(8).LineStyle = 1
$d1b.Item(10).ColorIndex = -4105
$d1b.Item(10).LineStyle = 1
$d1b.Item(9).ColorIndex = -4105
$d1b.Item(9).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").ClearFormats()
$c1b2 = $ws2.Range("C1").Borders
$c1b2.Item(8).ColorIndex = -4105
$c1b2.Item(8).LineStyle = 1
$c1b2.Item(9).ColorIndex = -4105
$c1b2.Item(9).LineStyle = 1

$ws2.Range("D1").ClearFormats()
$d1b2 = $ws2.Range("D1").Borders
$d1b2.Item(8).ColorIndex = -4105
$d1b2.Item(8).LineStyle = 1
$d1b2.Item(10).ColorIndex = -4105
$d1b2.Item(10).LineStyle = 1
$d1b2.Item(9).ColorIndex = -4105
$d1b2.Item(9).LineStyle = 1

$ws2.Range("F1").ClearFormats()
$f1b2 = $ws2.Range("F1").Borders
$f1b2.Item(8).ColorIndex = -4105
$f1b2.Item(8).LineStyle = 1
$f1b2.Item(9).ColorIndex = -4105
$f1b2.Item(9).LineStyle = 1

$ws2.Range("G1").ClearFormats()
$g1b2 = $ws2.Range("G1").Borders
$g1b2.Item(8).ColorIndex = -4105
$g1b2.Item(8).LineStyle = 1
$g1b2.Item(10).ColorIndex = -4105
$g1b2.Item(10).LineStyle = 1
$g1b2.Item(9).ColorIndex = -4105
$g1b2.Item(9).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
